# This presentation ships two DrawingML themes:
#   ppt/theme/theme1.xml -> the deck's real theme ("Integral"), wired to the
#                            slide master (and therefore every slide)
#   ppt/theme/theme2.xml -> the notes master's theme ("Office Theme")
#
# The target edit swaps the two themes' contents, so the slide master (and
# all slides) now use the stock "Office Theme" color palette instead of
# "Integral". Concretely only the 12-color <a:clrScheme> differs between the
# two themes (fonts/effects are already identical), so we push the "Office
# Theme" palette onto every slide's ThemeColorScheme, in clrScheme order:
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink

$p = $ppt.ActivePresentation

# Target palette ("Office Theme"), as 0xRRGGBB values, in clrScheme order.
$targetColorsRGB = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

# Apply across every slide so the shared theme is updated regardless of
# which slide/master the host keys the theme part off of.
$slideRange = $p.Slides.Range()
$themeColors = $slideRange.ThemeColorScheme

for ($i = 0; $i -lt $targetColorsRGB.Count; $i++) {
    $hex = $targetColorsRGB[$i]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    # PowerPoint's RGB color integers are packed as 0xBBGGRR.
    $vbaRgb = ($b -shl 16) -bor ($g -shl 8) -bor $r
    $themeColors.Item($i + 1).RGB = $vbaRgb
}
